$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("展览")
$ws4 = $wb.Worksheets.Item("全部类型")

$ws1.Range("F3").Value = 587
$ws1.Range("F5").Value = 301
$ws1.Range("F6").Value = 1118
$ws1.Range("F7").Value = 1462
$ws1.Range("F9").Value = 119
$ws1.Range("F10").Value = 761
$ws1.Range("F11").Value = 72
$ws1.Range("F12").Value = 182
$ws1.Range("F13").Value = 119
$ws1.Range("F14").Value = 461
$ws1.Range("F15").Value = 1406
$ws1.Range("F17").Value = 122
$ws1.Range("F18").Value = 284
$ws1.Range("F19").Value = 5211
$ws1.Range("F20").Value = 83
$ws1.Range("F21").Value = 666
$ws1.Range("F22").Value = 1019
$ws1.Range("F23").Value = 41
$ws1.Range("F24").Value = 256
$ws1.Range("F26").Value = 6048
$ws1.Range("F27").Value = 74
$ws1.Range("F28").Value = 127
$ws1.Range("F31").Value = 14754
$ws1.Range("F32").Value = 1467
$ws1.Range("F33").Value = 240
$ws1.Range("F34").Value = 110
$ws1.Range("F36").Value = 10389
$ws1.Range("F37").Value = 663
$ws1.Range("F38").Value = 4233
$ws1.Range("F39").Value = 168

$ws4.Range("F3").Value = 587
$ws4.Range("F5").Value = 301
$ws4.Range("F6").Value = 1118
$ws4.Range("F7").Value = 1462
$ws4.Range("F9").Value = 119
$ws4.Range("F10").Value = 761
$ws4.Range("F11").Value = 72
$ws4.Range("F12").Value = 182
$ws4.Range("F13").Value = 119
$ws4.Range("F14").Value = 461
$ws4.Range("F15").Value = 1406
$ws4.Range("F17").Value = 122
$ws4.Range("F18").Value = 284
$ws4.Range("F20").Value = 5211
$ws4.Range("F21").Value = 83
$ws4.Range("F22").Value = 666
$ws4.Range("F24").Value = 1019
$ws4.Range("F25").Value = 41
$ws4.Range("F26").Value = 256
$ws4.Range("F29").Value = 6048
$ws4.Range("F30").Value = 74
$ws4.Range("F31").Value = 127
$ws4.Range("F34").Value = 14754
$ws4.Range("F35").Value = 1467
$ws4.Range("F36").Value = 240
$ws4.Range("F37").Value = 110
$ws4.Range("F39").Value = 10393
$ws4.Range("F40").Value = 663
$ws4.Range("F41").Value = 4233
$ws4.Range("F42").Value = 168
